# Refresh the legacy GSC "Video Indexing" export.
# The underlying date window advanced by two days: the oldest two rows
# (2025-09-28 and 2025-09-29) fall out of range, every other row shifts
# up by two, and the two newest days (now rows 2 and 3) don't have
# indexing results back from Search Console yet, so their
# "No video indexed" / "Video indexed" columns come back blank - exactly
# like row 2 used to look before the refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Chart")

# Dropping rows 2:3 shifts all the later rows up by two and shrinks the
# used range from A1:D90 down to A1:D88 automatically.
$ws.Rows("2:3").Delete()

# The two newest days of data (now in rows 2 and 3) haven't been
# crawled/indexed yet, so those columns come back empty (an explicit
# empty string, matching the existing blank-data convention used by
# this export for not-yet-indexed days).
$ws.Range("B2").Text = ""
$ws.Range("C2").Text = ""
$ws.Range("B3").Text = ""
$ws.Range("C3").Text = ""
